# Update Active_Outages.xlsx per commit "Update Active_Outages.xlsx - 6/16/2025, 10:48:54 AM"
$wb = $excel.ActiveWorkbook

# ---- Sheet R1 ----
$ws1 = $wb.Worksheets.Item("R1")

# G2: 3875.6 -> 3876.0  (keep stored as text so the trailing .0 survives)
$ws1.Range("G2").NumberFormat = "@"
$ws1.Range("G2").Value = "3876.0"

# G3: 15.2 -> 15.5
$ws1.Range("G3").NumberFormat = "@"
$ws1.Range("G3").Value = "15.5"

# Row 4 gets populated with new outage data
$ws1.Range("B4").Value = "R4"
$ws1.Range("D4").Value = "asq0342"
$ws1.Range("I4").Value = "SCECO"
$ws1.Range("J4").Value = "In progress"
$ws1.Range("L4").Value = "Latis"

# ---- Sheet R6 ----
$ws6 = $wb.Worksheets.Item("R6")

# G2: 17.6 -> 17.4
$ws6.Range("G2").NumberFormat = "@"
$ws6.Range("G2").Value = "17.4"

# Row 3 gets cleared out (outage resolved/removed)
$ws6.Range("B3").Value = ""
$ws6.Range("D3").Value = ""
$ws6.Range("I3").Value = ""
$ws6.Range("J3").Value = ""
$ws6.Range("L3").Value = ""
